$d = $word.ActiveDocument

# 1) "{ now() }" -> "{now}" (header date formatting syntax)
$d.Content.Find.Execute("{ now() }", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{now}", 2)

# 2) "{month}" / "{year}" placeholders that were split across runs with
#    spell-check proofErr wrapping collapse to a single run's worth of text.
#    Find/Replace on the visible text normalizes this to "{month}" / "{year}".
$d.Content.Find.Execute("{month}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{month}", 2)
$d.Content.Find.Execute("{year}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{year}", 2)

# 3) "{LeadAuditor}" placeholder likewise collapses from three runs to one.
$d.Content.Find.Execute("{LeadAuditor}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{LeadAuditor}", 2)
